# Trace_Report_NOVAMEAL_initial.xlsx - refreshed trace-search result row.
# The car's last reported event moved from a Departure at OTTAWA, KS / train
# HKCKDE / destination LOVELAND to a "Placed Actual" event at JOHNSTOWN, CO
# (no train ID reported for this event), with an updated search timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_format_trace")

# Row 1: search-result description banner (date/time of the trace refreshed)
$ws.Range("A1").Value = "Description unknown, completed 10/16/2023 07:44:06 EDT, by WPJTOWN1.The search returned: 1 events."

# Row 3: the single data row
$ws.Range("C3").Value = "JOHNSTOWN"   # Location City: OTTAWA -> JOHNSTOWN
$ws.Range("D3").Value = "CO"          # State: KS -> CO
$ws.Range("F3").Value = 11            # Day: 3 -> 11
$ws.Range("G3").Value = 1511          # Time: 804 -> 1511
$ws.Range("H3").Value = "Placed Actual" # Event: Departure -> Placed Actual
$ws.Range("I3").ClearContents()       # Train ID: HKCKDE removed (none reported)
$ws.Range("J3").Value = "LOVELAND"    # Destination City: unchanged
$ws.Range("K3").Value = "CO"          # (State) unchanged
